# Fix descriptions for case 30 ("A30"/"L30") and case 31 ("A31"/"L31") on the
# 'descriptions' worksheet: the long-description text for these two rows had
# been entered swapped/mismatched with the age vs. length survey wording.
#
#   row 4  (letter/number "A30", lookup key "A30") -> "fish & surv length"
#   row 5  (letter/number "A31", lookup key "A31") -> "fish & .5surv length"
#   row 19 (letter/number "L30", lookup key "L30") -> "fish & surv length"
#   row 20 (letter/number "L31", lookup key "L31") -> "fish & .5surv length"
#
# All of the 'scenarios' worksheet's LOOKUP(...) formulas that reference
# descriptions!$D:$D recalc automatically once these source cells change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("descriptions")

$ws.Range("D4").Value  = "fish & surv length"
$ws.Range("D5").Value  = "fish & .5surv length"
$ws.Range("D19").Value = "fish & surv length"
$ws.Range("D20").Value = "fish & .5surv length"
